# Automatic update of files.
#
# The underlying data rows (2, 5-14) got re-synced from source and ended up
# re-ordered: the content that used to live in one row now lives in another
# row (rows 3 and 4 are unaffected). A couple of cells also changed value as
# part of that re-sync (Q/R coordinates, and the "Publik kommentar" (AC) free
# text note moved along with its row).
#
# Strategy: snapshot the values that actually vary across the affected rows
# (columns A, B, D, E, F, G, H, Q, R, AC) from every source row first, then
# write them back into their new row positions. Reading everything up front
# avoids clobbering source data before it has been captured, since this is a
# full permutation with no fixed points among rows 2 / 5-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSnapshot($row) {
    return @{
        A  = $ws.Range("A$row").Value2
        B  = $ws.Range("B$row").Value2
        D  = $ws.Range("D$row").Value2
        E  = $ws.Range("E$row").Value2
        F  = $ws.Range("F$row").Value2
        G  = $ws.Range("G$row").Value2
        H  = $ws.Range("H$row").Value2
        Q  = $ws.Range("Q$row").Value2
        R  = $ws.Range("R$row").Value2
        AC = $ws.Range("AC$row").Value2
    }
}

# Capture the "before" content of every affected row.
$snapshots = @{}
foreach ($row in 2, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14) {
    $snapshots[$row] = Get-RowSnapshot $row
}

# old row -> new row mapping (rows 3 and 4 are untouched).
$mapping = @{
    2  = 9
    5  = 8
    6  = 10
    7  = 11
    8  = 5
    9  = 14
    10 = 13
    11 = 12
    12 = 6
    13 = 2
    14 = 7
}

function Set-RowFromSnapshot($row, $data) {
    $ws.Range("A$row").Value = $data.A
    $ws.Range("B$row").Value = $data.B
    $ws.Range("D$row").Value = $data.D
    $ws.Range("E$row").Value = $data.E
    $ws.Range("F$row").Value = $data.F
    $ws.Range("G$row").Value = $data.G
    $ws.Range("H$row").Value = $data.H
    $ws.Range("Q$row").Value = $data.Q
    $ws.Range("R$row").Value = $data.R
    if ($data.AC) {
        $ws.Range("AC$row").Value = $data.AC
    } else {
        $ws.Range("AC$row").Value = ""
    }
}

foreach ($oldRow in $mapping.Keys) {
    $newRow = $mapping[$oldRow]
    Set-RowFromSnapshot $newRow $snapshots[$oldRow]
}
